# Regenerate save_data: replace the old "Strike#" derived K column (column G)
# values with newly-computed K values (rows 2-74), and update the two
# dependent cells (H36, I36) that also changed as part of the regen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (column G) value
$kValues = @{
    2  = 2
    3  = 1
    4  = 3
    5  = 1
    6  = 3
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 3
    15 = 1
    16 = 0
    17 = 3
    18 = 3
    19 = 2
    20 = 1
    21 = 2
    22 = 2
    23 = 3
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 2
    30 = 2
    31 = 2
    32 = 3
    33 = 0
    34 = 3
    35 = 2
    36 = 8
    37 = 3
    38 = 2
    39 = 3
    40 = 2
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 2
    48 = 0
    49 = 1
    50 = 0
    51 = 0
    52 = 0
    53 = 4
    54 = 1
    55 = 2
    56 = 2
    57 = 4
    58 = 0
    59 = 0
    60 = 1
    61 = 2
    62 = 1
    63 = 0
    64 = 5
    65 = 1
    66 = 0
    67 = 2
    68 = 0
    69 = 1
    70 = 1
    71 = 1
    72 = 2
    73 = 0
    74 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

# Row 36 also needed its IP (H) and I0 (I) values regenerated alongside K.
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 5
